$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 becomes "line7" (was "extr1"): from_bus/to_bus values updated
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 becomes "line8" (was "extr2"): from_bus/to_bus values updated, now in service
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10 stays "extr1" but from_bus/to_bus/in_service change
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 stays "extr2" but from_bus/to_bus/in_service change
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 stays "extr3" but from_bus/in_service change
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

# Row 13 stays "extr4" but to_bus/in_service change
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14 stays "extr5" but from_bus/to_bus/in_service change
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15 stays "extr6" but from_bus/to_bus/in_service change
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New row 16: "extr7"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true
$ws.Range("A10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 17: "extr8"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
$ws.Range("A10").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "done"
